$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.274.35'
$ws.Range("E2").Value = '  -1.93%  '

# Row 3
$ws.Range("D3").Value = '2.433.15'
$ws.Range("E3").Value = '  -1.55%  '

# Row 4
$ws.Range("E4").Value = '  -0.29%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.13'
$ws.Range("E5").Value = '  -2.26%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.04'
$ws.Range("E6").Value = '  -3.45%  '

# Row 7
$ws.Range("E7").Value = '  +0.10%  '

# Row 8
$ws.Range("E8").Value = '  -2.51%  '

# Row 9
$ws.Range("D9").Value = '2.428.84'
$ws.Range("E9").Value = '  -2.06%  '

# Row 10
$ws.Range("E10").Value = '  -5.54%  '

# Row 11
$ws.Range("E11").Value = '  +0.98%  '

# Row 12
$ws.Range("B12").Value = 'Toncoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.18'
$ws.Range("E12").Value = '  -2.94%  '

# Row 13
$ws.Range("B13").Value = 'Cardano'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.352'
$ws.Range("E13").Value = '  -3.11%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.47'
$ws.Range("E14").Value = '  -2.84%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000173'
$ws.Range("E15").Value = '  -6.37%  '

# Row 16
$ws.Range("D16").Value = '2.871.21'
$ws.Range("E16").Value = '  -1.56%  '

# Row 17
$ws.Range("D17").Value = '62.160.28'
$ws.Range("E17").Value = '  -1.88%  '

# Row 18
$ws.Range("D18").Value = '2.426.07'
$ws.Range("E18").Value = '  -1.94%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.04'
$ws.Range("E19").Value = '  -4.52%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.09'
$ws.Range("E20").Value = '  -2.86%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.92'
$ws.Range("E21").Value = '  -1.50%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.12'
$ws.Range("E22").Value = '  -2.94%  '

# Row 23
$ws.Range("E23").Value = '  +3.89%  '

# Row 24
$ws.Range("E24").Value = '  +0.22%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.92'
$ws.Range("E25").Value = '  -3.74%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '621.55'
$ws.Range("E26").Value = '  -1.07%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.90'
$ws.Range("E27").Value = '  +1.73%  '

# Row 28
$ws.Range("D28").Value = '2.549.68'
$ws.Range("E28").Value = '  -1.77%  '

# Row 29
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0953'
$ws.Range("E29").Value = '  -9.66%  '

# Row 30
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.38%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.44'
$ws.Range("E31").Value = '  -5.75%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.02'
$ws.Range("E32").Value = '  -5.18%  '

# Row 33
$ws.Range("E33").Value = '  -3.62%  '

# Row 34
$ws.Range("E34").Value = '  -7.53%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.00'
$ws.Range("E35").Value = '  -3.66%  '

# Row 36
$ws.Range("E36").Value = '  +0.14%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.45'
$ws.Range("E37").Value = '  -7.27%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.375'
$ws.Range("E38").Value = '  -2.90%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.65'
$ws.Range("E39").Value = '  -1.64%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '146.57'
$ws.Range("E40").Value = '  -0.40%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.21'
$ws.Range("E41").Value = '  -5.67%  '

# Row 42
$ws.Range("E42").Value = '  -6.75%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.59'
$ws.Range("E43").Value = '  +1.79%  '

# Row 44
$ws.Range("E44").Value = '  +0.01%  '

# Row 45
$ws.Range("E45").Value = '  -8.50%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '145.22'
$ws.Range("E46").Value = '  -3.56%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.68'
$ws.Range("E47").Value = '  -2.38%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0521'
$ws.Range("E48").Value = '  -5.53%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.96'
$ws.Range("E49").Value = '  -5.89%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.592'
$ws.Range("E50").Value = '  -2.94%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0228'
$ws.Range("E51").Value = '  -5.26%  '
